$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row at the end (row 63) with the "Sting 20 rs" item (Drinks category)
# first, so that its strings land in the shared-string table ahead of "TAKA TAK".
$ws.Cells.Item(63, 1).Value = "Drinks"
$ws.Cells.Item(63, 2).Value = "Sting 20 rs"
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 20
$ws.Cells.Item(63, 5).Value = "Sting 20 rs.jpg"
$ws.Cells.Item(63, 6).Value = "Fast Food"

# Insert a new row above the current row 50 ("Melody toffee 1 rs"),
# shifting existing rows 50-63 down to 51-64.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the "TAKA TAK" item (Snacks category).
$ws.Cells.Item(50, 1).Value = "Snacks"
$ws.Cells.Item(50, 2).Value = "TAKA TAK"
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 10
$ws.Cells.Item(50, 5).Value = "TAKA TAK.jpg"
$ws.Cells.Item(50, 6).Value = "Fast Food"

# The hidden "_FilterDatabase" name grows by one row because of the row
# inserted above (it does not include the appended row at the very end).
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "=Sheet1!`$A`$1:`$F`$61"

# Update the view to match the target state: scrolled down to the new rows,
# with the active selection on the newly inserted row's last cell.
$excel.Goto($ws.Range("A49"), $true)
$ws.Range("F50").Select()
